$d = $word.ActiveDocument
$d.Paragraphs.Item(2).Range.Delete()
